# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# D/E columns hold text-formatted numbers ("75.894.32", "  +1.73%  ", etc.),
# so NumberFormat is forced to Text ('@') before writing each D-cell to stop
# Excel auto-coercing the string into a float (which would e.g. turn
# "1.00" into 1 or "0.200" into 0.2). Style is reset to 'Normal' right after
# so the cell doesn't end up carrying a lingering Text-format style index.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '75.894.32'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.910.65'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.20%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '198.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '597.85'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.71%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.200'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '2.914.26'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.435'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +18.46%  '
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('E13').Value = '  +2.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.448.54'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '75.787.79'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000192'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.44'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.918.62'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.94'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.76'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '379.38'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.30'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.26'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.24'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.064.80'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.20'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.65'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000108'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.69%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('E31').Value = '  +3.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '503.92'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.71'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.57%  '
$ws.Range('E34').Value = '  +2.41%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '164.84'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.26'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.56%  '
$ws.Range('B38').Value = 'Cronos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.106'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +24.79%  '
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.66'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.93%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.113'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.71%  '
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '179.97'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.80%  '
$ws.Range('E43').Value = '  +2.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.99'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.49%  '
$ws.Range('E45').Value = '  -0.58%  '
$ws.Range('E46').Value = '  +1.44%  '
$ws.Range('E47').Value = '  -1.32%  '
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.570'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.662'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.72'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.54%  '
